# Update cryptocurrency price (D) and volume-change (E) columns with the
# latest scraped values. Column D holds plain text-like numbers (e.g.
# "25.988.21") that must remain text, so the range is temporarily switched
# to text format while assigning, then restored to the default style so the
# cells end up unstyled exactly like the rest of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "25.994.66"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.639.59"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").Value = "214.53"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").Value = "0.5087"
$ws.Range("E6").Value = "  +0.75%  "
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").Value = "0.2556"
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("D9").Value = "0.06338"
$ws.Range("E9").Value = "  -0.81%  "
$ws.Range("D10").Value = "19.55"
$ws.Range("E10").Value = "  -0.31%  "
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("D12").Value = "4.269"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").Value = "1.643.90"
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("D14").Value = "0.5415"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").Value = "63.99"
$ws.Range("E15").Value = "  -1.26%  "
$ws.Range("D16").Value = "0.0₅7672"
$ws.Range("E16").Value = "  -2.40%  "
$ws.Range("D17").Value = "26.005.29"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("D19").Value = "198.24"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "4.411"
$ws.Range("E20").Value = "  +0.41%  "
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("D22").Value = "6.035"
$ws.Range("E22").Value = "  +0.82%  "
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").Value = "1.864"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").Value = "141.07"
$ws.Range("E25").Value = "  +0.80%  "
$ws.Range("D26").Value = "0.1187"
$ws.Range("E26").Value = "  +3.88%  "
$ws.Range("D27").Value = "6.799"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("D28").Value = "15.60"
$ws.Range("E28").Value = "  -0.78%  "
$ws.Range("E29").Value = "  -0.71%  "
$ws.Range("D30").Value = "0.04886"
$ws.Range("E30").Value = "  -0.98%  "
$ws.Range("D31").Value = "3.249"
$ws.Range("E31").Value = "  -0.59%  "
$ws.Range("D32").Value = "3.162"
$ws.Range("E32").Value = "  -1.10%  "
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").Value = "0.9039"
$ws.Range("E35").Value = "  +1.15%  "
$ws.Range("E36").Value = "  -0.87%  "
$ws.Range("D37").Value = "1.142.75"
$ws.Range("E37").Value = "  +0.33%  "
$ws.Range("D38").Value = "0.5430"
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("D39").Value = "0.01559"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").Value = "2.523"
$ws.Range("E41").Value = "  -1.65%  "
$ws.Range("D42").Value = "0.0₈128"
$ws.Range("E42").Value = "  +6.47%  "
$ws.Range("D43").Value = "0.8087"
$ws.Range("E43").Value = "  -1.26%  "
$ws.Range("D44").Value = "99.10"
$ws.Range("E44").Value = "  -0.26%  "
$ws.Range("E45").Value = "  -4.80%  "
$ws.Range("D46").Value = "1.778.85"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").Value = "54.73"
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("D50").Value = "0.05111"
$ws.Range("E50").Value = "  +1.13%  "
$ws.Range("D51").Value = "1.001"
$ws.Range("E51").Value = "  -0.65%  "

$priceRange.Style = "Normal"
